$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 573.80554
$ws.Range("J17").Value = 309.17856
$ws.Range("L17").Value = 927.53568
$ws.Range("N17").Value = -1263.53568
$ws.Range("H21").Value = 9729.083000000001
$ws.Range("I21").Value = 8031.1875
$ws.Range("J21").Value = 13124.875
$ws.Range("K21").Value = 8031.1875
$ws.Range("L21").Value = 13124.875
$ws.Range("M21").Value = -7563.1875
$ws.Range("N21").Value = -14060.875
$ws.Range("H23").Value = 9729.083000000001
$ws.Range("I23").Value = 8031.1875
$ws.Range("J23").Value = 13124.875
$ws.Range("K23").Value = 8031.1875
$ws.Range("L23").Value = 13124.875
$ws.Range("M23").Value = -7797.1875
$ws.Range("N23").Value = -13592.875
$ws.Range("H62").Value = 1728
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 1728
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H98").Value = 4821.6514
$ws.Range("I98").Value = 2123.5925
$ws.Range("J98").Value = 9374.625
$ws.Range("K98").Value = 2123.5925
$ws.Range("L98").Value = 9374.625
$ws.Range("M98").Value = -625.5925000000002
$ws.Range("N98").Value = -12370.625
$ws.Range("H112").Value = 1307.6852
$ws.Range("J112").Value = 1322.9246
$ws.Range("L112").Value = 3968.7738
$ws.Range("N112").Value = -6184.7738
$ws.Range("H122").Value = 4821.6514
$ws.Range("I122").Value = 2123.5925
$ws.Range("J122").Value = 9374.625
$ws.Range("K122").Value = 6370.7775
$ws.Range("L122").Value = 28123.875
$ws.Range("M122").Value = -3920.7775
$ws.Range("N122").Value = -33023.875
$ws.Range("H123").Value = 41223.332
$ws.Range("J123").Value = 41835
$ws.Range("L123").Value = 41835
$ws.Range("N123").Value = -51635
$ws.Range("H127").Value = 1741.6
$ws.Range("J127").Value = 2290
$ws.Range("L127").Value = 6870
$ws.Range("N127").Value = -16790
$ws.Range("H129").Value = 836.84
$ws.Range("J129").Value = 865.12634
$ws.Range("L129").Value = 2595.37902
$ws.Range("N129").Value = -12595.37902
$ws.Range("H137").Value = 1490263.2
$ws.Range("I137").Value = 1765089.8
$ws.Range("J137").Value = 6200
$ws.Range("K137").Value = 5295269.4
$ws.Range("L137").Value = 18600
$ws.Range("M137").Value = -5292719.4
$ws.Range("N137").Value = -23700
$ws.Range("H141").Value = 36377.035
$ws.Range("I141").Value = 41530.16
$ws.Range("K141").Value = 124590.48
$ws.Range("M141").Value = -119410.48

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2320.5715
$ws.Range("I45").Value = 2886
$ws.Range("J45").Value = 907
$ws.Range("K45").Value = 2886
$ws.Range("L45").Value = 907
$ws.Range("M45").Value = -2509
$ws.Range("N45").Value = -1661
$ws.Range("H97").Value = 2003.4
$ws.Range("I97").Value = 1002.875
$ws.Range("K97").Value = 1002.875
$ws.Range("M97").Value = -506.875

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2999.8333
$ws.Range("J107").Value = 2999.75
$ws.Range("L107").Value = 2999.75
$ws.Range("N107").Value = -6839.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5907.206
$ws.Range("I31").Value = 1127.45
$ws.Range("J31").Value = 12735.429
$ws.Range("K31").Value = 1127.45
$ws.Range("L31").Value = 12735.429
$ws.Range("M31").Value = -832.45
$ws.Range("N31").Value = -13325.429
$ws.Range("H34").Value = 5907.206
$ws.Range("I34").Value = 1127.45
$ws.Range("J34").Value = 12735.429
$ws.Range("K34").Value = 1127.45
$ws.Range("L34").Value = 12735.429
$ws.Range("M34").Value = -925.45
$ws.Range("N34").Value = -13139.429
$ws.Range("H53").Value = 57842
$ws.Range("J53").Value = 57842
$ws.Range("L53").Value = 57842
$ws.Range("N53").Value = -59056
$ws.Range("H122").Value = 3493.125
$ws.Range("I122").Value = 1889
$ws.Range("J122").Value = 6166.6665
$ws.Range("K122").Value = 5667
$ws.Range("L122").Value = 18499.9995
$ws.Range("M122").Value = -3217
$ws.Range("N122").Value = -23399.9995
$ws.Range("H132").Value = 2683.0952
$ws.Range("I132").Value = 1709.1875
$ws.Range("J132").Value = 5799.6
$ws.Range("K132").Value = 5127.5625
$ws.Range("L132").Value = 17398.8
$ws.Range("M132").Value = -2597.5625
$ws.Range("N132").Value = -22458.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 9166.333000000001
$ws.Range("J80").Value = 8749.5
$ws.Range("L80").Value = 26248.5
$ws.Range("N80").Value = -28120.5
$ws.Range("H83").Value = 9166.333000000001
$ws.Range("J83").Value = 8749.5
$ws.Range("L83").Value = 78745.5
$ws.Range("N83").Value = -88105.5
$ws.Range("H113").Value = 3572146.8
$ws.Range("I113").Value = 589.2273
$ws.Range("J113").Value = 9616321
$ws.Range("K113").Value = 1767.6819
$ws.Range("L113").Value = 28848963
$ws.Range("M113").Value = 402.3181
$ws.Range("N113").Value = -28853303
$ws.Range("H122").Value = 3116.262
$ws.Range("J122").Value = 3452.3057
$ws.Range("L122").Value = 31070.7513
$ws.Range("N122").Value = -35970.7513
$ws.Range("H125").Value = 2000
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 6000
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -1080
$ws.Range("N125").ClearContents()
$ws.Range("H131").Value = 795.3299
$ws.Range("I131").Value = 316.1111
$ws.Range("J131").Value = 844.3409
$ws.Range("K131").Value = 948.3333
$ws.Range("L131").Value = 2533.0227
$ws.Range("M131").Value = 4091.6667
$ws.Range("N131").Value = -12613.0227

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 6583.4
$ws.Range("I33").Value = 3000
$ws.Range("J33").Value = 7479.25
$ws.Range("K33").Value = 3000
$ws.Range("L33").Value = 7479.25
$ws.Range("M33").Value = -2748
$ws.Range("N33").Value = -7983.25
$ws.Range("H132").Value = 4396.24
$ws.Range("I132").Value = 3757.524
$ws.Range("J132").Value = 7749.5
$ws.Range("K132").Value = 11272.572
$ws.Range("L132").Value = 23248.5
$ws.Range("M132").Value = -8742.572
$ws.Range("N132").Value = -28308.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 1877.5555
$ws.Range("J22").Value = 2122.4443
$ws.Range("K22").Value = 1877.5555
$ws.Range("L22").Value = 2122.4443
$ws.Range("M22").Value = -1582.5555
$ws.Range("N22").Value = -2712.4443
$ws.Range("H27").Value = 2000
$ws.Range("I27").Value = 1877.5555
$ws.Range("J27").Value = 2122.4443
$ws.Range("K27").Value = 1877.5555
$ws.Range("L27").Value = 2122.4443
$ws.Range("M27").Value = -1770.5555
$ws.Range("N27").Value = -2336.4443
$ws.Range("H82").Value = 6308.3335
$ws.Range("I82").Value = 7949.643
$ws.Range("J82").Value = 3025.7144
$ws.Range("K82").Value = 7949.643
$ws.Range("L82").Value = 3025.7144
$ws.Range("M82").Value = -7588.643
$ws.Range("N82").Value = -3747.7144
$ws.Range("H85").Value = 6308.3335
$ws.Range("I85").Value = 7949.643
$ws.Range("J85").Value = 3025.7144
$ws.Range("K85").Value = 7949.643
$ws.Range("L85").Value = 3025.7144
$ws.Range("M85").Value = -6701.643
$ws.Range("N85").Value = -5521.7144

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2496.1904
$ws.Range("I81").Value = 2266.5881
$ws.Range("J81").Value = 3472
$ws.Range("K81").Value = 4533.1762
$ws.Range("L81").Value = 6944
$ws.Range("M81").Value = -3472.1762
$ws.Range("N81").Value = -9066
$ws.Range("H84").Value = 2496.1904
$ws.Range("I84").Value = 2266.5881
$ws.Range("J84").Value = 3472
$ws.Range("K84").Value = 22665.881
$ws.Range("L84").Value = 34720
$ws.Range("M84").Value = -17361.881
$ws.Range("N84").Value = -45328
$ws.Range("H122").Value = 8589.929
$ws.Range("I122").Value = 7350.5
$ws.Range("J122").Value = 10242.5
$ws.Range("K122").Value = 22051.5
$ws.Range("L122").Value = 30727.5
$ws.Range("M122").Value = -19601.5
$ws.Range("N122").Value = -35627.5
